$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.906.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.24%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.68%  "

# Row 4
$ws.Range("E4").Value = "  -1.98%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5050"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.006"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.90%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2577"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.16%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06404"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.50%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07792"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.54%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.276"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.76%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.642.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.49%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.867.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.67%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5441"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.14%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7925"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.40%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.31%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.984.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.98%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.86%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "197.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.06%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.378"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.93%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.915"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.46%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.973"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.61%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.28%  "

# Row 25
$ws.Range("E25").Value = "  -4.47%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1137"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.36%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.841"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.58%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.243"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04948"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.41%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.269"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.27%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.203"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "

# Row 34
$ws.Range("E34").Value = "  +0.27%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.376"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.04%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8938"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.22%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.614"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.31%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.146.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.32%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5556"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.68%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01569"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.71%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.007"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.07%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.710"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.23%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8124"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.23%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.46%  "

# Row 45
$ws.Range("E45").Value = "  +7.24%  "

# Row 46
$ws.Range("E46").Value = "  +0.66%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4534"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.23%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9998"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.32%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.14%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05067"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.40%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.006"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.63%  "
